# "Gestion du Projet.xlsx" - update home page progress tracking table.
#
# - Module News (article, categorie)   D13 : 0%  -> 60%  (now "in progress" -> Neutre style)
# - Module Don (tableau des dons)      D14 : 0%  -> 20%  (now "in progress" -> Neutre style)
# - Mise en page                      D22 : (blank) -> 50% (now "in progress" -> Neutre style)
# - Selection cursor left on H8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D9 already carries the "in progress" look (percentage number format + the
# "Neutre" yellow cell style). Clone that exact formatting onto the three
# cells whose progress moved off 0%, instead of re-applying the named style
# (which would not reinstate the percentage number format), so they end up
# sharing the very same style record rather than minting new ones.
$ws.Range("D9").Copy()
$ws.Range("D13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New progress values.
$ws.Range("D13").Value = 0.6
$ws.Range("D14").Value = 0.2
$ws.Range("D22").Value = 0.5

# Leave the active selection on H8, as in the saved file.
$ws.Range("H8").Select()
